# Integrate Agora data for elec sector variable BGDPbES (BAU Guaranteed
# Dispatch Percentage by Electricity Source).
#
# The Agora-sourced update removes the previously "guaranteed" 100%
# dispatch assumption for nuclear and hydro: both sources now have 0%
# guaranteed dispatch in the BAU case (2015, column B) which then ripples
# through the shared "=copy across the row" formulas in columns C:AK via
# recalculation.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BGDPbES")

# Row 4 = nuclear, Row 5 = hydro. Only the base year (column B) needs to be
# written -- the rest of the row (C:AK) are shared formulas that simply
# reference column B and recalc automatically.
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# The refreshed workbook also turns on iterative calculation.
$excel.Iteration = $true
$excel.MaxChange = 0.00001

# Restore the selections left behind in each sheet, ending on the "About"
# sheet so it remains the active tab.
$ws.Activate()
$ws.Range("B6").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A9").Select()
